$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Tumor" query text in B3 (SamplesTab row) so that the tumor
# status is read directly from the sample node instead of the collected
# list that was built earlier in the query.
$newTumorQuery = @'
MATCH (s:study)<--(p:participant)<--(samp:sample)
WHERE s.study_name in ["Washington University PDX Development and Trial Center"]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@

$ws.Range("B3").Value = $newTumorQuery

# Move the active selection, matching the saved cursor position in the
# workbook (B12).
$ws.Range("B12").Select()
